$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, whether an apostrophe prefix
# is needed to force Excel to keep a numeric-looking string as text.
$updates = @(
    ,@('D2', '54.359.92', $false)
    ,@('E2', '  -3.15%  ', $false)
    ,@('D3', '2.247.00', $false)
    ,@('E3', '  -4.76%  ', $false)
    ,@('D4', '1.01', $true)
    ,@('E4', '  +0.59%  ', $false)
    ,@('D5', '494.95', $true)
    ,@('E5', '  -1.96%  ', $false)
    ,@('D6', '127.37', $true)
    ,@('E6', '  -3.67%  ', $false)
    ,@('D7', '0.998', $true)
    ,@('E7', '  -0.06%  ', $false)
    ,@('D8', '0.530', $true)
    ,@('E8', '  -1.89%  ', $false)
    ,@('D9', '2.298.18', $false)
    ,@('E9', '  -3.66%  ', $false)
    ,@('E10', '  -0.97%  ', $false)
    ,@('E12', '  +0.84%  ', $false)
    ,@('D13', '4.63', $true)
    ,@('E13', '  -3.97%  ', $false)
    ,@('D14', '2.658.53', $false)
    ,@('E14', '  -4.51%  ', $false)
    ,@('D15', '21.70', $true)
    ,@('E15', '  -0.46%  ', $false)
    ,@('D16', '54.272.88', $false)
    ,@('E16', '  -3.26%  ', $false)
    ,@('D17', '0.0000129', $true)
    ,@('E17', '  -2.09%  ', $false)
    ,@('D18', '2.289.10', $false)
    ,@('E18', '  -3.03%  ', $false)
    ,@('D19', '10.00', $true)
    ,@('E19', '  +0.51%  ', $false)
    ,@('E20', '  +0.49%  ', $false)
    ,@('D21', '304.61', $true)
    ,@('E21', '  -1.51%  ', $false)
    ,@('D22', '6.48', $true)
    ,@('E22', '  +3.73%  ', $false)
    ,@('D23', '0.995', $true)
    ,@('E23', '  -0.49%  ', $false)
    ,@('E24', '  -2.67%  ', $false)
    ,@('D25', '62.95', $true)
    ,@('E25', '  -3.72%  ', $false)
    ,@('E26', '  +0.66%  ', $false)
    ,@('E27', '  -0.06%  ', $false)
    ,@('E28', '  +2.80%  ', $false)
    ,@('D29', '2.398.07', $false)
    ,@('E29', '  -3.72%  ', $false)
    ,@('E30', '  -1.86%  ', $false)
    ,@('D31', '170.78', $true)
    ,@('E31', '  +0.24%  ', $false)
    ,@('E32', '  -2.60%  ', $false)
    ,@('D33', '0.0₃0687', $false)
    ,@('E33', '  -3.75%  ', $false)
    ,@('E34', '  +1.05%  ', $false)
    ,@('D36', '0.993', $true)
    ,@('E36', '  -0.52%  ', $false)
    ,@('E37', '  -1.49%  ', $false)
    ,@('D38', '17.60', $true)
    ,@('E38', '  -0.57%  ', $false)
    ,@('D39', '1.21', $true)
    ,@('E39', '  +1.55%  ', $false)
    ,@('D40', '0.868', $true)
    ,@('E40', '  -0.78%  ', $false)
    ,@('E41', '  -2.01%  ', $false)
    ,@('D42', '35.56', $true)
    ,@('E42', '  -2.47%  ', $false)
    ,@('D43', '0.376', $true)
    ,@('E43', '  -0.33%  ', $false)
    ,@('E44', '  -1.27%  ', $false)
    ,@('D45', '129.45', $true)
    ,@('E45', '  +2.63%  ', $false)
    ,@('E46', '  -1.05%  ', $false)
    ,@('D47', '4.90', $true)
    ,@('E47', '  -0.84%  ', $false)
    ,@('E48', '  -0.17%  ', $false)
    ,@('D49', '0.550', $true)
    ,@('E49', '  -2.05%  ', $false)
    ,@('D50', '242.55', $true)
    ,@('E50', '  -0.58%  ', $false)
    ,@('E51', '  -0.38%  ', $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    if ($forceText) {
        $val = "'" + $val
    }
    $ws.Range($ref).Value = $val
}
